$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75 and Row 76 swap their product details (description, code, category, brand),
# reflecting the re-ordering of entries, while quantity/turnover are updated to new totals.

$ws.Range("A75").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("E75").Value = "Παπαγάλος® Ελληνικός Καφές Κουπάτος 143gr"
$ws.Range("F75").Value = "5201219486417"
$ws.Range("G75").Value = 3.7
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 50
$ws.Range("J75").Value = "Παπαγάλος"
$ws.Range("K75").Value = 18
$ws.Range("L75").Value = 30.32

$ws.Range("A76").Value = "Πελάτες Τιμή Πώλησης"
$ws.Range("E76").Value = "Papadopoulou® Cookies με κομματια σοκολατας 180gr"
$ws.Range("F76").Value = "5201004021755"
$ws.Range("G76").Value = 1.69
$ws.Range("H76").Value = 1
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = "Papadopoulou"
$ws.Range("K76").Value = 20
$ws.Range("L76").Value = 17.6

# Row 79 brand reference updates (text value unchanged, stays "Papadopoulou")
$ws.Range("J79").Value = "Papadopoulou"

# Updated cumulative figures further down the sheet
$ws.Range("K81").Value = 59.896
$ws.Range("L81").Value = 69.06

$ws.Range("K82").Value = 84.43000000000001
$ws.Range("L82").Value = 77.68000000000001

$ws.Range("K83").Value = 415.276
$ws.Range("L83").Value = 773.8600000000001
